$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I (9th column) so the existing
# "comments" / "questions that need to be asked to the PI" column
# shifts from I to J, and a new "new_value_label" column takes over
# column I.
$ws.Range("I1").EntireColumn.Insert()

# Match the formatting of the neighboring column (H) for the new
# header and description cells before filling in their text
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("H3").Copy()
$ws.Range("I3").PasteSpecial(-4122)  # xlPasteFormats

# New header text in the (now empty) column I
$ws.Range("I1").Value = "new_value_label"

# New description text in row 3 for column I
$ws.Range("I3").Value = "If novalue labels are provided, give a label.provide the value and the label, seperated by comma, similar to the value_label column"

# Column widths: column I keeps (approximately) the same width as
# column H (17.1796875 characters), and the new column J (previously
# holding the "comments"/questions text) gets a wider width of 25
$ws.Range("I1").ColumnWidth = 16.9
$ws.Range("J1").ColumnWidth = 24.15

# Adjust row 3 height to fit the new content
$ws.Range("A3").RowHeight = 116

# Clear the clipboard marching-ants left over from the copy operations
$excel.CutCopyMode = $false

# Update the active selection to I3, matching the authored change
$ws.Range("I3").Select()
